$wb = $excel.ActiveWorkbook

# --- Remove leftover empty placeholder cells B2/B3 on "ODI Batting" ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B3").ClearContents()

# --- Add the new "ODI Batting Extra" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Headers
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Copy the existing header formatting (bold font + thin border + centered
# alignment) from another sheet's header row so the new header row matches
# the workbook's established header style exactly.
$headerStyleSource = $odiBatting.Range("A1:F1")
$headerStyleSource.Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# Row 2 -> match 4669
$extra.Range("A2").Value = "'4669"
$extra.Range("B2").Value = 10
$extra.Range("C2").Value = "'"
$extra.Range("D2").Value = "'"
$extra.Range("E2").Value = "'"
$extra.Range("F2").Value = "NO"

# Row 3 -> match 4673
$extra.Range("A3").Value = "'4673"
$extra.Range("B3").Value = "'"
$extra.Range("C3").Value = "'"
$extra.Range("D3").Value = "'"
$extra.Range("E3").Value = "'"
$extra.Range("F3").Value = "NO"

# Row 4 -> match 4676
$extra.Range("A4").Value = "'4676"
$extra.Range("B4").Value = 10
$extra.Range("C4").Value = "'0"
$extra.Range("D4").Value = "'1"
$extra.Range("E4").Value = "'4.11%"
$extra.Range("F4").Value = "NO"

# The leading-apostrophe text entries above flag their cells with a
# "number stored as text" quote-prefix style. Strip that back off (while
# keeping the cells as Text) by pasting the plain default format from an
# untouched cell on top of them - mirrors what a user would do via
# "Clear > Formats" after typing text-as-text in real Excel.
$plainFormat = $extra.Range("H10")
$plainFormat.Copy()
$extra.Range("A2:A4").PasteSpecial(-4122)
$extra.Range("B3").PasteSpecial(-4122)
$extra.Range("C2:E4").PasteSpecial(-4122)
